# CQ_InputTemplate.xlsx update
# The underlying medium names "Medium-A" (columns B:G) and "Medium-B"
# (columns H:M) in the Media-name block (rows 33-40) are consolidated into
# a single medium called "MED-A" across the whole B33:M40 range.
# Everything else (AB1_AB2_AB3 labels, concentration labels, and the
# CONCATENATE summary formulas in rows 57-64) are derived/unchanged and
# will recompute automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B33:M40").Value = "MED-A"

# Restore the view to where the user ended up editing (bottom of sheet,
# near the summary table), matching the saved selection in the workbook.
$ws.Range("O38").Select()
